# "updated 4.0 files and mdl"
#
# - About!C1: bump the "last updated" date serial 45320 -> 45392
# - MCF: several plant-type max capacity factors bumped from 0.85/0.95 -> 1
#   (B2,B3,B4,B6,B10,B11,B12,B13,B14,B16,B17,B18). The dependent formula
#   cells (B19,B20,B21,B22,B24,B25 = B2/B4/B10/B14) recalc automatically.
# - MCF sheet: active selection moves from E8 to B17

$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$mcf   = $wb.Worksheets.Item("MCF")

# About sheet: updated date stamp
$about.Range("C1").Value = 45392

# MCF sheet: raise these max capacity factors to 1 (100%)
$mcf.Range("B2").Value  = 1   # hard coal
$mcf.Range("B3").Value  = 1   # natural gas steam turbine
$mcf.Range("B4").Value  = 1   # natural gas combined cycle
$mcf.Range("B6").Value  = 1   # hydro
$mcf.Range("B10").Value = 1   # biomass
$mcf.Range("B11").Value = 1   # geothermal
$mcf.Range("B12").Value = 1   # petroleum
$mcf.Range("B13").Value = 1   # natural gas peaker
$mcf.Range("B14").Value = 1   # lignite
$mcf.Range("B16").Value = 1   # crude oil
$mcf.Range("B17").Value = 1   # heavy or residual fuel oil
$mcf.Range("B18").Value = 1   # municipal solid waste

# Move/save the active selection on the MCF sheet to B17
$mcf.Activate()
$mcf.Range("B17").Select()
